$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 533
$ws.Range("I28").Value = 533
$ws.Range("K28").Value = 533
$ws.Range("M28").Value = -48
$ws.Range("H40").Value = 1949.4595
$ws.Range("I40").Value = 1926.875
$ws.Range("J40").Value = 1991.1538
$ws.Range("K40").Value = 1926.875
$ws.Range("L40").Value = 1991.1538
$ws.Range("M40").Value = -1751.875
$ws.Range("N40").Value = -2341.1538
$ws.Range("H64").Value = 3994.2856
$ws.Range("I64").Value = 4240
$ws.Range("J64").Value = 3666.6667
$ws.Range("K64").Value = 4240
$ws.Range("L64").Value = 3666.6667
$ws.Range("M64").Value = -3992
$ws.Range("N64").Value = -4162.6667
$ws.Range("H67").Value = 3994.2856
$ws.Range("I67").Value = 4240
$ws.Range("J67").Value = 3666.6667
$ws.Range("K67").Value = 4240
$ws.Range("L67").Value = 3666.6667
$ws.Range("M67").Value = -3382
$ws.Range("N67").Value = -5382.6667
$ws.Range("H76").Value = 2796.1304
$ws.Range("I76").Value = 2747.9048
$ws.Range("K76").Value = 2747.9048
$ws.Range("M76").Value = -2432.9048
$ws.Range("H79").Value = 2796.1304
$ws.Range("I79").Value = 2747.9048
$ws.Range("K79").Value = 2747.9048
$ws.Range("M79").Value = -1655.9048
$ws.Range("H125").Value = 1278.3914
$ws.Range("I125").Value = 805.1429000000001
$ws.Range("J125").Value = 2014.5555
$ws.Range("K125").Value = 7246.2861
$ws.Range("L125").Value = 18130.9995
$ws.Range("M125").Value = -4786.2861
$ws.Range("N125").Value = -23050.9995
$ws.Range("H132").Value = 4548180
$ws.Range("I132").Value = 5407525
$ws.Range("K132").Value = 16222575
$ws.Range("M132").Value = -16220045
$ws.Range("H137").Value = 3808.3416
$ws.Range("I137").Value = 4441.4614
$ws.Range("K137").Value = 13324.3842
$ws.Range("M137").Value = -10774.3842

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4145.66
$ws.Range("I32").Value = 4145.66
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4145.66
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3858.66
$ws.Range("N32").ClearContents()
$ws.Range("H45").Value = 1655.7931
$ws.Range("I45").Value = 1240.16
$ws.Range("J45").Value = 4253.5
$ws.Range("K45").Value = 1240.16
$ws.Range("L45").Value = 4253.5
$ws.Range("M45").Value = -863.1600000000001
$ws.Range("N45").Value = -5007.5
$ws.Range("H46").Value = 6500
$ws.Range("J46").Value = 6500
$ws.Range("L46").Value = 6500
$ws.Range("N46").Value = -7138
$ws.Range("H61").Value = 2626.074
$ws.Range("I61").Value = 1918.5555
$ws.Range("J61").Value = 4041.111
$ws.Range("K61").Value = 1918.5555
$ws.Range("L61").Value = 4041.111
$ws.Range("M61").Value = -1706.5555
$ws.Range("N61").Value = -4465.111
$ws.Range("H63").Value = 2499.8
$ws.Range("I63").Value = 2124.75
$ws.Range("K63").Value = 2124.75
$ws.Range("M63").Value = -1438.75
$ws.Range("H66").Value = 2499.8
$ws.Range("I66").Value = 2124.75
$ws.Range("K66").Value = 10623.75
$ws.Range("M66").Value = -7191.75
$ws.Range("H88").Value = 3243.2856
$ws.Range("I88").Value = 2099
$ws.Range("J88").Value = 4769
$ws.Range("K88").Value = 2099
$ws.Range("L88").Value = 4769
$ws.Range("M88").Value = -1693
$ws.Range("N88").Value = -5581
$ws.Range("H91").Value = 3243.2856
$ws.Range("I91").Value = 2099
$ws.Range("J91").Value = 4769
$ws.Range("K91").Value = 2099
$ws.Range("L91").Value = 4769
$ws.Range("M91").Value = -695
$ws.Range("N91").Value = -7577
$ws.Range("H97").Value = 474.64
$ws.Range("I97").Value = 429.25
$ws.Range("J97").Value = 656.2
$ws.Range("K97").Value = 429.25
$ws.Range("L97").Value = 656.2
$ws.Range("M97").Value = 66.75
$ws.Range("N97").Value = -1648.2
$ws.Range("H136").Value = 2626.074
$ws.Range("I136").Value = 1918.5555
$ws.Range("J136").Value = 4041.111
$ws.Range("K136").Value = 5755.666499999999
$ws.Range("L136").Value = 12123.333
$ws.Range("M136").Value = -3205.666499999999
$ws.Range("N136").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 1619.6666
$ws.Range("J17").Value = 2004.5
$ws.Range("L17").Value = 2004.5
$ws.Range("N17").Value = -2348.5
$ws.Range("H80").Value = 665.7083
$ws.Range("J80").Value = 653.4286
$ws.Range("L80").Value = 653.4286
$ws.Range("N80").Value = -2649.4286
$ws.Range("H83").Value = 665.7083
$ws.Range("J83").Value = 653.4286
$ws.Range("L83").Value = 3267.143
$ws.Range("N83").Value = -13251.143

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1741.1237
$ws.Range("I31").Value = 986.12964
$ws.Range("J31").Value = 2689.2559
$ws.Range("K31").Value = 986.12964
$ws.Range("L31").Value = 2689.2559
$ws.Range("M31").Value = -691.12964
$ws.Range("N31").Value = -3279.2559
$ws.Range("H34").Value = 1741.1237
$ws.Range("I34").Value = 986.12964
$ws.Range("J34").Value = 2689.2559
$ws.Range("K34").Value = 986.12964
$ws.Range("L34").Value = 2689.2559
$ws.Range("M34").Value = -784.12964
$ws.Range("N34").Value = -3093.2559
$ws.Range("H99").Value = 4701.4
$ws.Range("I99").Value = 2049.75
$ws.Range("K99").Value = 2049.75
$ws.Range("M99").Value = -551.75
$ws.Range("H126").Value = 4701.4
$ws.Range("I126").Value = 2049.75
$ws.Range("K126").Value = 6149.25
$ws.Range("M126").Value = -3679.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4471.8667
$ws.Range("I70").Value = 4516.4443
$ws.Range("J70").Value = 4405
$ws.Range("K70").Value = 4516.4443
$ws.Range("L70").Value = 4405
$ws.Range("M70").Value = -4246.4443
$ws.Range("N70").Value = -4945
$ws.Range("H73").Value = 4471.8667
$ws.Range("I73").Value = 4516.4443
$ws.Range("J73").Value = 4405
$ws.Range("K73").Value = 4516.4443
$ws.Range("L73").Value = 4405
$ws.Range("M73").Value = -3580.4443
$ws.Range("N73").Value = -6277
$ws.Range("H80").Value = 3131.0789
$ws.Range("I80").Value = 2687
$ws.Range("J80").Value = 3985.077
$ws.Range("K80").Value = 2687
$ws.Range("L80").Value = 3985.077
$ws.Range("M80").Value = -1689
$ws.Range("N80").Value = -5981.077
$ws.Range("H83").Value = 3131.0789
$ws.Range("I83").Value = 2687
$ws.Range("J83").Value = 3985.077
$ws.Range("K83").Value = 13435
$ws.Range("L83").Value = 19925.385
$ws.Range("M83").Value = -8443
$ws.Range("N83").Value = -29909.385
$ws.Range("H132").Value = 3521.2727
$ws.Range("I132").Value = 3390.6667
$ws.Range("K132").Value = 10172.0001
$ws.Range("M132").Value = -7642.000100000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 35000
$ws.Range("J87").Value = 35000
$ws.Range("L87").Value = 35000
$ws.Range("N87").Value = -37246
$ws.Range("H90").Value = 35000
$ws.Range("J90").Value = 35000
$ws.Range("L90").Value = 105000
$ws.Range("N90").Value = -116232
$ws.Range("H132").Value = 3964.95
$ws.Range("I132").Value = 3144.3333
$ws.Range("J132").Value = 4636.364
$ws.Range("K132").Value = 9432.999899999999
$ws.Range("L132").Value = 13909.092
$ws.Range("M132").Value = -6902.999899999999
$ws.Range("N132").Value = -18969.092

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1873.3334
$ws.Range("I81").Value = 1540
$ws.Range("J81").Value = 2111.4285
$ws.Range("K81").Value = 3080
$ws.Range("L81").Value = 4222.857
$ws.Range("M81").Value = -2019
$ws.Range("N81").Value = -6344.857
$ws.Range("H84").Value = 1873.3334
$ws.Range("I84").Value = 1540
$ws.Range("J84").Value = 2111.4285
$ws.Range("K84").Value = 15400
$ws.Range("L84").Value = 21114.285
$ws.Range("M84").Value = -10096
$ws.Range("N84").Value = -31722.285
